# Planning.xlsx edit: "Removed builds and start of untangling report"
#
# Adds a new "Hekkenberg" meeting note row to the Quotes sheet (with the
# Dutch/English "assume the reader is stupid..." reminder quote plus the
# meeting date it came from), widens column A on that sheet to fit the new
# author name, and switches the active/selected tab from "Mensen in de
# loop" over to "Quotes".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quotes")

# New row of data - enter in Author, Quote(english), Bron, Quote(dutch)
# order so new shared-string entries line up the same way the original
# author's edit produced them.
$ws.Range("A4").Value = "Hekkenberg"
$ws.Range("C4").Value = "Assume the reader is stupid, lazy and oblivious"
$ws.Range("D4").Value = "Meeting 11-04-2018"
$ws.Range("B4").Value = "Ga ervan uit dat de lezer dom, lui en vergeetachtig is"

# Widen column A so the new author name isn't clipped.
$ws.Columns.Item(1).ColumnWidth = 20.6

# Put the selection where the author left it and make Quotes the active
# (selected) sheet in the workbook - this moves tabSelected from "Mensen
# in de loop" to "Quotes" and updates the workbook's activeTab.
$ws.Range("B2").Select() | Out-Null
$ws.Activate() | Out-Null
